# Update the Vtn-Itga5 NATMI interaction table with newly computed TPM-based statistics,
# and drop the obsolete "Resolving-Mac" sending-cluster block (rows that become empty /
# no longer produced once the upstream script was re-run with the new TPM values).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: ECs -> ECs
$ws.Range("G2").Value = 7.844453333333334
$ws.Range("H2").Value = 23.53336
$ws.Range("I2").Value = 0.1489290605659587
$ws.Range("J2").Value = 0.1489290605659588
$ws.Range("M2").Value = 35.04689966666667
$ws.Range("N2").Value = 105.140699
$ws.Range("O2").Value = 0.3824629895491901
$ws.Range("P2").Value = 0.3824629895491901
$ws.Range("Q2").Value = 274.9237689131822
$ws.Range("R2").Value = 2474.31392021864
$ws.Range("S2").Value = 0.05695985373480898
$ws.Range("T2").Value = 0.05695985373480898

# Row 3: ECs -> FAPs
$ws.Range("G3").Value = 7.844453333333334
$ws.Range("H3").Value = 23.53336
$ws.Range("I3").Value = 0.1489290605659587
$ws.Range("J3").Value = 0.1489290605659588
$ws.Range("M3").Value = 29.913269
$ws.Range("N3").Value = 89.739807
$ws.Range("O3").Value = 0.3264402385872224
$ws.Range("P3").Value = 0.3264402385872223
$ws.Range("Q3").Value = 234.6532427179467
$ws.Range("R3").Value = 2111.87918446152
$ws.Range("S3").Value = 0.04861643806372246
$ws.Range("T3").Value = 0.04861643806372246

# Row 4: ECs -> MuSCs
$ws.Range("G4").Value = 7.844453333333334
$ws.Range("H4").Value = 23.53336
$ws.Range("I4").Value = 0.1489290605659587
$ws.Range("J4").Value = 0.1489290605659588
$ws.Range("M4").Value = 8.911727666666666
$ws.Range("N4").Value = 26.735183
$ws.Range("O4").Value = 0.09725271102035077
$ws.Range("P4").Value = 0.09725271102035075
$ws.Range("Q4").Value = 69.90763180054222
$ws.Range("R4").Value = 629.16868620488
$ws.Range("S4").Value = 0.0144837548897535
$ws.Range("T4").Value = 0.0144837548897535

# Row 5: ECs -> Resolving-Mac
$ws.Range("G5").Value = 7.844453333333334
$ws.Range("H5").Value = 23.53336
$ws.Range("I5").Value = 0.1489290605659587
$ws.Range("J5").Value = 0.1489290605659588
$ws.Range("M5").Value = 17.76285166666667
$ws.Range("N5").Value = 53.288555
$ws.Range("O5").Value = 0.1938440608432367
$ws.Range("P5").Value = 0.1938440608432367
$ws.Range("Q5").Value = 139.3398609660889
$ws.Range("R5").Value = 1254.0587486948
$ws.Range("S5").Value = 0.02886901387767379
$ws.Range("T5").Value = 0.02886901387767379

# Row 6: FAPs -> ECs
$ws.Range("G6").Value = 20.35396833333334
$ws.Range("H6").Value = 61.06190500000001
$ws.Range("I6").Value = 0.3864255740794268
$ws.Range("J6").Value = 0.3864255740794268
$ws.Range("M6").Value = 35.04689966666667
$ws.Range("N6").Value = 105.140699
$ws.Range("O6").Value = 0.3824629895491901
$ws.Range("P6").Value = 0.3824629895491901
$ws.Range("Q6").Value = 713.343485996844
$ws.Range("R6").Value = 6420.091373971597
$ws.Range("S6").Value = 0.1477934803006796
$ws.Range("T6").Value = 0.1477934803006796

# Row 7: FAPs -> FAPs
$ws.Range("G7").Value = 20.35396833333334
$ws.Range("H7").Value = 61.06190500000001
$ws.Range("I7").Value = 0.3864255740794268
$ws.Range("J7").Value = 0.3864255740794268
$ws.Range("M7").Value = 29.913269
$ws.Range("N7").Value = 89.739807
$ws.Range("O7").Value = 0.3264402385872224
$ws.Range("P7").Value = 0.3264402385872223
$ws.Range("Q7").Value = 608.8537299724818
$ws.Range("R7").Value = 5479.683569752336
$ws.Range("S7").Value = 0.1261448565986925
$ws.Range("T7").Value = 0.1261448565986924

# Row 8: FAPs -> MuSCs
$ws.Range("G8").Value = 20.35396833333334
$ws.Range("H8").Value = 61.06190500000001
$ws.Range("I8").Value = 0.3864255740794268
$ws.Range("J8").Value = 0.3864255740794268
$ws.Range("M8").Value = 8.911727666666666
$ws.Range("N8").Value = 26.735183
$ws.Range("O8").Value = 0.09725271102035077
$ws.Range("P8").Value = 0.09725271102035075
$ws.Range("Q8").Value = 181.3890227226239
$ws.Range("R8").Value = 1632.501204503615
$ws.Range("S8").Value = 0.03758093468681964
$ws.Range("T8").Value = 0.03758093468681964

# Row 9: FAPs -> Resolving-Mac
$ws.Range("G9").Value = 20.35396833333334
$ws.Range("H9").Value = 61.06190500000001
$ws.Range("I9").Value = 0.3864255740794268
$ws.Range("J9").Value = 0.3864255740794268
$ws.Range("M9").Value = 17.76285166666667
$ws.Range("N9").Value = 53.288555
$ws.Range("O9").Value = 0.1938440608432367
$ws.Range("P9").Value = 0.1938440608432367
$ws.Range("Q9").Value = 361.5445203330306
$ws.Range("R9").Value = 3253.900682997275
$ws.Range("S9").Value = 0.0749063024932351
$ws.Range("T9").Value = 0.07490630249323509

# Row 10: MuSCs -> ECs
$ws.Range("G10").Value = 24.47399366666667
$ws.Range("H10").Value = 73.421981
$ws.Range("I10").Value = 0.4646453653546145
$ws.Range("J10").Value = 0.4646453653546145
$ws.Range("M10").Value = 35.04689966666667
$ws.Range("N10").Value = 105.140699
$ws.Range("O10").Value = 0.3824629895491901
$ws.Range("P10").Value = 0.3824629895491901
$ws.Range("Q10").Value = 857.7376004783022
$ws.Range("R10").Value = 7719.63840430472
$ws.Range("S10").Value = 0.1777096555137015
$ws.Range("T10").Value = 0.1777096555137015

# Row 11: MuSCs -> FAPs
$ws.Range("G11").Value = 24.47399366666667
$ws.Range("H11").Value = 73.421981
$ws.Range("I11").Value = 0.4646453653546145
$ws.Range("J11").Value = 0.4646453653546145
$ws.Range("M11").Value = 29.913269
$ws.Range("N11").Value = 89.739807
$ws.Range("O11").Value = 0.3264402385872224
$ws.Range("P11").Value = 0.3264402385872223
$ws.Range("Q11").Value = 732.0971560552964
$ws.Range("R11").Value = 6588.874404497667
$ws.Range("S11").Value = 0.1516789439248075
$ws.Range("T11").Value = 0.1516789439248075

# Row 12: MuSCs -> MuSCs
$ws.Range("G12").Value = 24.47399366666667
$ws.Range("H12").Value = 73.421981
$ws.Range("I12").Value = 0.4646453653546145
$ws.Range("J12").Value = 0.4646453653546145
$ws.Range("M12").Value = 8.911727666666666
$ws.Range("N12").Value = 26.735183
$ws.Range("O12").Value = 0.09725271102035077
$ws.Range("P12").Value = 0.09725271102035075
$ws.Range("Q12").Value = 218.1055664730581
$ws.Range("R12").Value = 1962.950098257523
$ws.Range("S12").Value = 0.04518802144377763
$ws.Range("T12").Value = 0.04518802144377762

# Row 13: MuSCs -> Resolving-Mac
$ws.Range("G13").Value = 24.47399366666667
$ws.Range("H13").Value = 73.421981
$ws.Range("I13").Value = 0.4646453653546145
$ws.Range("J13").Value = 0.4646453653546145
$ws.Range("M13").Value = 17.76285166666667
$ws.Range("N13").Value = 53.288555
$ws.Range("O13").Value = 0.1938440608432367
$ws.Range("P13").Value = 0.1938440608432367
$ws.Range("Q13").Value = 434.7279191919395
$ws.Range("R13").Value = 3912.551272727455
$ws.Range("S13").Value = 0.09006874447232785
$ws.Range("T13").Value = 0.09006874447232785

# The source data for the "Resolving-Mac" sending cluster (previously rows 14-17)
# is no longer present after the script re-run, so remove those rows entirely.
$ws.Rows("14:17").Delete()

Write-Host "Edit complete"